# chore: adapt column header formatting to respective input file names (#7)
#
# - Rename the "_old" / "_new" header-name suffixes in row 1 to the concrete
#   format-version identifiers "_FV2404" / "_FV2410".
# - Turn the data range A1:U69 into a native Excel Table ("Table1") so the
#   (now renamed) header row drives the table's column names.
# - Freeze the header row (split/freeze below row 1) and keep the selection
#   anchored in the lower-left frozen pane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) ------------------------------------
# Columns A-J carried the "_old" suffix, column K is the plain "diff" column
# (left untouched), and columns L-U carried the "_new" suffix.
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# --- 2. Turn the range into a proper Excel Table ("Table1") ----------------
# Must happen after the header rename above so the ListObject picks up the
# new column names straight from the header row.
$tableRange = $ws.Range("A1:U69")
$listObj = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObj.Name = "Table1"

# --- 3. Freeze the header row -----------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
